$d = $word.ActiveDocument

# Locate the end of the paragraph that finishes with
# "...plusieurs profils en même temps. " (the paragraph that needs to be
# split in two). Using Find keeps this robust against any paragraph
# numbering/index assumptions.
$anchor = $d.Content
$found = $anchor.Find.Execute("plusieurs profils en même temps. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the end of the target paragraph"
}
$splitPoint = $anchor.End

# Split the paragraph in two by inserting a paragraph break right after
# the found text (and before the existing paragraph mark).
$breakRange = $d.Range($splitPoint, $splitPoint)
$breakRange.InsertParagraphAfter()

# The new (currently empty) paragraph now sits right after the split
# point; fill it in with the new text, then drop the bookmark in the
# middle of "voulait" exactly as in the target revision.
$newParaStart = $splitPoint + 1
$newPara = $d.Range($newParaStart, $newParaStart)
$part1 = "Nous avons décidé d’utiliser des doubles cliques pour pouvoir afficher le détail des personnages car l’utilisateur peut sans le faire exprès cliquer sur un profil alors qu’il voul"
$part2 = "ait pas voir son détail."
$newPara.InsertAfter($part1 + $part2)

$bookmarkPos = $newParaStart + $part1.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
